$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("devices")

# Add the new values required by the diff (these also register new shared strings)
$ws.Range("B2").Value = "Android"
$ws.Range("K2").Value = "Appium"
$ws.Range("L2").Value = "testssts"

# Move the active selection to L3, matching the author's final cursor position
$ws.Activate()
$ws.Range("L3").Select()
